$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = 333
$ws.Range("B4").Value = "Vlad"
